$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.155.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.457.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.85"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.90%  "
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  +2.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.046.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("E13").Value = "  +2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.452.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.142.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "74.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.590.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000116"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.182"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "31.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.491.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0771"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.794"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.583.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.68%  "
$ws.Range("E51").Value = "  +0.06%  "
